$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.324.97'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '3.514.01'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.72%  '
$ws.Range('E7').Value = '  -0.60%  '
$ws.Range('D8').Value = '3.507.00'
$ws.Range('E8').Value = '  -0.51%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('E10').Value = '  -1.46%  '
$ws.Range('E11').Value = '  +8.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.583'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '46.21'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000275'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.44%  '
$ws.Range('D15').Value = '4.086.22'
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.30'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '611.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.05%  '
$ws.Range('D18').Value = '3.520.45'
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('D19').Value = '70.424.35'
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('E20').Value = '  +0.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('E22').Value = '  -1.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.06'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -19.57%  '
$ws.Range('E24').Value = '  -1.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '97.20'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('E26').Value = '  -4.30%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('E28').Value = '  -3.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.12'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.99'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.14'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.99'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '647.26'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +13.30%  '
$ws.Range('E34').Value = '  -0.89%  '
$ws.Range('E35').Value = '  -4.74%  '
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0996'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.71'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0476'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '56.66'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('E42').Value = '  +1.52%  '
$ws.Range('D43').Value = '0.0₃0744'
$ws.Range('E43').Value = '  +4.97%  '
$ws.Range('D44').Value = '3.364.61'
$ws.Range('E44').Value = '  -0.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.310'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.91'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '32.19'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.55'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.129'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '134.20'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.01%  '
